# Rewrites the TOC bookmark names (w:name="_Toc...") to match a
# refreshed "Table of Contents" numbering, as produced by Word when it
# regenerates the hidden _Toc bookmarks. The bookmark ids and the
# ranges they cover are left untouched; only the names change.
#
# Word's COM object model does not support assigning Bookmark.Name
# directly (it is read-only in practice), so we recreate each bookmark
# at the same Start/End with the new name, which preserves its id and
# position in the XML.

$d = $word.ActiveDocument

$renameMap = @(
    @{ Old = "_Toc5728486"; New = "_Toc5791565" },
    @{ Old = "_Toc5728487"; New = "_Toc5791566" },
    @{ Old = "_Toc5728488"; New = "_Toc5791567" },
    @{ Old = "_Toc5728489"; New = "_Toc5791568" },
    @{ Old = "_Toc5728490"; New = "_Toc5791569" },
    @{ Old = "_Toc5728491"; New = "_Toc5791570" },
    @{ Old = "_Toc5728492"; New = "_Toc5791571" },
    @{ Old = "_Toc5728493"; New = "_Toc5791572" },
    @{ Old = "_Toc5728494"; New = "_Toc5791573" },
    @{ Old = "_Toc5728495"; New = "_Toc5791574" },
    @{ Old = "_Toc5728496"; New = "_Toc5791575" },
    @{ Old = "_Toc5728497"; New = "_Toc5791576" },
    @{ Old = "_Toc5728498"; New = "_Toc5791577" },
    @{ Old = "_Toc5728499"; New = "_Toc5791578" },
    @{ Old = "_Toc5728500"; New = "_Toc5791579" },
    @{ Old = "_Toc5728501"; New = "_Toc5791580" },
    @{ Old = "_Toc5728502"; New = "_Toc5791581" },
    @{ Old = "_Toc5728503"; New = "_Toc5791582" },
    @{ Old = "_Toc5728504"; New = "_Toc5791583" },
    @{ Old = "_Toc5728505"; New = "_Toc5791584" }
)

foreach ($pair in $renameMap) {
    $oldName = $pair.Old
    $newName = $pair.New

    if ($d.Bookmarks.Exists($oldName)) {
        $bm = $d.Bookmarks.Item($oldName)
        $start = $bm.Start
        $end = $bm.End
        $bm.Delete()
        $newRange = $d.Range($start, $end)
        $d.Bookmarks.Add($newName, $newRange)
        Write-Output "Renamed $oldName -> $newName"
    } else {
        Write-Output "MISSING $oldName"
    }
}
